$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =====================================================================
# New crime data collected - weekly CompStat report refresh
#  - Bump the report "Volume/Number" and the covered week's date range
#  - Refresh every crime-category statistic in the Week-to-Date /
#    28-Day / Year-to-Date table (rows 15-27) to the newly collected
#    figures, including derived percent-change columns
# =====================================================================

# ---- Masthead: Volume 30 Number 17 -> Number 18 ----
$ws.Range("A8").Characters(21, 2).Text = "18"

# ---- Masthead: reporting week 4/24/2023-4/30/2023 -> 5/1/2023-5/7/2023 ----
# (replace the second date first so the first replacement's offset,
#  which comes earlier in the string, is not shifted by a length change)
$ws.Range("C9").Characters(47, 9).Text = "5/7/2023"
$ws.Range("C9").Characters(27, 9).Text = "5/1/2023"

# ---- Four cells that flip between the numeric "0/***.* " placeholder
#      and a real number (or vice versa). A plain .Value assignment of
#      a numeric-looking string is auto-coerced to a real number, which
#      would leave the cell's number format (style) wrong, so each of
#      these pastes the format (only) from a same-column neighbour that
#      already carries the desired style after writing the raw value. ----

# F15: numeric 2 -> text placeholder "0" (style must stay the General/
# right-top style used by its row-mate C15, which is already that text)
$ws.Range("F15").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("F15").PasteSpecial(-4122)

# D18: text placeholder "0" -> numeric 10 (style must become the plain
# number style used by its column neighbours D17/D19)
$ws.Range("D18").Value = 10
$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)

# E18: text placeholder "***.*" -> numeric -80 (style must become the
# percent-change style used by its column neighbours E17/E19)
$ws.Range("E18").Value = -80
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)

# C26: text placeholder "0" -> numeric 1 (style must become the plain
# number style used by its column neighbours C25/C27)
$ws.Range("C26").Value = 1
$ws.Range("C25").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Bulk numeric cell updates (crime stats table, rows 15-27) ----
# Row 15 - Rape
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -37.5

# Row 16 - Robbery
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 42
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = -19.230769230769
$ws.Range("L16").Value = 44.827586206896
$ws.Range("M16").Value = 82.608695652173
$ws.Range("N16").Value = -86.229508196721

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 43
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 2.380952380952
$ws.Range("L17").Value = 2.380952380952
$ws.Range("M17").Value = 115
$ws.Range("N17").Value = -27.118644067796

# Row 18 - Burglary
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 101
$ws.Range("K18").Value = -43.564356435643
$ws.Range("L18").Value = 46.153846153846
$ws.Range("M18").Value = -13.636363636363
$ws.Range("N18").Value = -81.311475409836

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 23
$ws.Range("E19").Value = -8
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 95
$ws.Range("H19").Value = -5.263157894736
$ws.Range("I19").Value = 387
$ws.Range("J19").Value = 398
$ws.Range("K19").Value = -2.763819095477
$ws.Range("L19").Value = 78.341013824884
$ws.Range("M19").Value = 1.842105263157
$ws.Range("N19").Value = -67.396798652064

# Row 20 - G.L.A.
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 23
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = 53.333333333333
$ws.Range("L20").Value = 109.090909090909
$ws.Range("M20").Value = 187.5
$ws.Range("N20").Value = -92.281879194630

# Row 21 - TOTAL
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -27.272727272727
$ws.Range("F21").Value = 132
$ws.Range("G21").Value = 147
$ws.Range("H21").Value = -10.204081632653
$ws.Range("I21").Value = 557
$ws.Range("J21").Value = 616
$ws.Range("K21").Value = -9.577922077922
$ws.Range("L21").Value = 61.449275362318
$ws.Range("M21").Value = 12.072434607645
$ws.Range("N21").Value = -74.272517321016

# Row 22 - Transit
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -12.5
$ws.Range("I22").Value = 26
$ws.Range("J22").Value = 37
$ws.Range("K22").Value = -29.729729729729
$ws.Range("L22").Value = 18.181818181818
$ws.Range("M22").Value = 8.333333333333

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 80
$ws.Range("D24").Value = 78
$ws.Range("E24").Value = 2.564102564102
$ws.Range("F24").Value = 292
$ws.Range("G24").Value = 331
$ws.Range("H24").Value = -11.782477341389
$ws.Range("I24").Value = 1303
$ws.Range("J24").Value = 1312
$ws.Range("K24").Value = -0.685975609756
$ws.Range("L24").Value = 101.391035548686
$ws.Range("M24").Value = 124.268502581756

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 8.571428571428
$ws.Range("I25").Value = 117
$ws.Range("J25").Value = 119
$ws.Range("K25").Value = -1.680672268907
$ws.Range("L25").Value = 17
$ws.Range("M25").Value = 37.647058823529

# Row 26 - UCR Rape*
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 7
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = -22.222222222222
$ws.Range("L26").Value = 0

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 12
$ws.Range("H27").Value = -7.692307692307
$ws.Range("I27").Value = 32
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -3.030303030303
$ws.Range("L27").Value = 28
